$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 22 ("Biaya Layanan Tokopedia"),
# shifting everything below down by two rows.
$ws.Rows("22:23").Insert()

# Copy the formatting (borders, alignment, number formats) of the row
# directly above the insertion point down into the two new rows so they
# match the rest of the bordered table.
$ws.Range("A21:F21").Copy()
$ws.Range("A22:F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new line items.
$ws.Range("C22").Value = "Wing Bolt"
$ws.Range("D22").Value = 21
$ws.Range("E22").Value = 17200

$ws.Range("C23").Value = "Rubber Pad"
$ws.Range("D23").Value = 22
$ws.Range("E23").Value = 13100

# "Biaya Layanan Tokopedia" (now shifted to row 24) price increased.
$ws.Range("E24").Value = 10000

# Update the view to match where the author left off editing.
$ws.Range("I22").Select()
$excel.ActiveWindow.ScrollRow = 16

Write-Host "done"
